# Trade #64 closed at 2026-02-17 21:11:28 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results.xlsx" workbook to reflect
# a trade close event (MarketMaking strategy) plus the opening of a brand
# new trade right afterwards. It touches four sheets:
#   1. Summary          - aggregate capital / P&L / trade counters
#   2. Strategy Status   - per-strategy summary row for MarketMaking
#   3. All Trades        - full trade log (closes trade row 93, appends row 126)
#   4. MarketMaking       - per-strategy trade log (closes trade row 60, appends row 93)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a date-looking / time-looking string into a cell WITHOUT
# letting Excel auto-convert it into a date/time serial number. The
# source workbook stores these as plain text, so we force the cell to
# Text format before assigning the value.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# =======================================================================
# 1) Summary sheet
# =======================================================================
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B3").Value = 1401.07   # Current Capital
$summary.Range("B4").Value = 0.87      # Total P&L $
$summary.Range("B5").Value = 0.19      # Total P&L %
$summary.Range("B6").Value = 92        # Total Trades
$summary.Range("B8").Value = 36        # Losing Trades
$summary.Range("B9").Value = 47.83     # Win Rate %

# =======================================================================
# 2) Strategy Status sheet - MarketMaking row (row 5)
# =======================================================================
$status = $wb.Worksheets.Item("Strategy Status")

$status.Range("C5").Value = 101.07
$status.Range("D5").Value = 59
$status.Range("E5").Value = 0.76
$status.Range("F5").Value = 1.07
$status.Range("G5").Value = 50.85

# =======================================================================
# 3) All Trades sheet
# =======================================================================
$allTrades = $wb.Worksheets.Item("All Trades")

# --- Update existing row 93 (Trade #92) to reflect the close ---
$allTrades.Cells.Item(93, 7).Value  = 0.87            # G93 Exit Price
$allTrades.Cells.Item(93, 8).Value  = "CLOSED"        # H93 Status
$allTrades.Cells.Item(93, 9).Value  = -4.3956         # I93 P&L %
$allTrades.Cells.Item(93, 10).Value = -0.04           # J93 P&L $
$allTrades.Cells.Item(93, 11).Value = 101.07          # K93 Capital After
$allTrades.Cells.Item(93, 12).Value = "early_exit"    # L93 Exit Reason
$allTrades.Cells.Item(93, 13).Value = 0.14            # M93 Duration (min)

# --- Append new row 126 (Trade #125, newly opened) ---
Set-TextValue $allTrades.Cells.Item(126, 2) "2026-02-17"
Set-TextValue $allTrades.Cells.Item(126, 3) "21:11:21"
$allTrades.Cells.Item(126, 1).Value  = 125             # A126 Trade #
$allTrades.Cells.Item(126, 4).Value  = "MarketMaking"  # D126 Strategy
$allTrades.Cells.Item(126, 5).Value  = "UP"            # E126 Side
$allTrades.Cells.Item(126, 6).Value  = 0.91            # F126 Entry Price
$allTrades.Cells.Item(126, 8).Value  = "OPEN"          # H126 Status
$allTrades.Cells.Item(126, 9).Value  = 0               # I126 P&L %
$allTrades.Cells.Item(126, 10).Value = 0               # J126 P&L $
$allTrades.Cells.Item(126, 11).Value = 101.1146450978375 # K126 Capital After
$allTrades.Cells.Item(126, 13).Value = 0               # M126 Duration (min)
$allTrades.Cells.Item(126, 14).Value = 0               # N126 Entry Slippage (bps)
$allTrades.Cells.Item(126, 15).Value = 0               # O126 Exit Slippage (bps)
$allTrades.Cells.Item(126, 16).Value = 0.6             # P126 Confidence
$allTrades.Cells.Item(126, 17).Value = "Normal spread capture: 19600 bps" # Q126 Entry Reason

# =======================================================================
# 4) MarketMaking sheet
# =======================================================================
$mm = $wb.Worksheets.Item("MarketMaking")

# --- Update existing row 60 (Trade #92) to reflect the close ---
$mm.Cells.Item(60, 7).Value  = 0.87            # G60 Exit Price
$mm.Cells.Item(60, 8).Value  = "CLOSED"        # H60 Status
$mm.Cells.Item(60, 9).Value  = -4.3956         # I60 P&L %
$mm.Cells.Item(60, 10).Value = -0.04           # J60 P&L $
$mm.Cells.Item(60, 11).Value = 101.07          # K60 Capital After
$mm.Cells.Item(60, 16).Value = "early_exit"    # P60 Exit Reason
$mm.Cells.Item(60, 17).Value = 0.14            # Q60 Duration (min)

# --- Append new row 93 (Trade #125, newly opened) ---
Set-TextValue $mm.Cells.Item(93, 2) "2026-02-17"
Set-TextValue $mm.Cells.Item(93, 3) "21:11:21"
$mm.Cells.Item(93, 1).Value  = 125              # A93 Trade #
$mm.Cells.Item(93, 4).Value  = "MarketMaking"   # D93 Strategy
$mm.Cells.Item(93, 5).Value  = "UP"             # E93 Side
$mm.Cells.Item(93, 6).Value  = 0.91             # F93 Entry Price
$mm.Cells.Item(93, 8).Value  = "OPEN"           # H93 Status
$mm.Cells.Item(93, 9).Value  = 0                # I93 P&L %
$mm.Cells.Item(93, 10).Value = 0                # J93 P&L $
$mm.Cells.Item(93, 11).Value = 101.1146450978375 # K93 Capital After
$mm.Cells.Item(93, 12).Value = 0                # L93 Entry Slippage (bps)
$mm.Cells.Item(93, 13).Value = 0                # M93 Exit Slippage (bps)
$mm.Cells.Item(93, 14).Value = 0.6              # N93 Confidence
$mm.Cells.Item(93, 15).Value = "Normal spread capture: 19600 bps" # O93 Entry Reason
$mm.Cells.Item(93, 17).Value = 0                # Q93 Duration (min)
